$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pvalue_Table_Ancestry")

# Updated p-value / fdr figures (recomputed downstream analysis)
$ws.Range("C2").Value = 0.99231324419116695
$ws.Range("D2").Value = 0.99231324419116695
$ws.Range("H2").Value = 0.77055369341804303
$ws.Range("I2").Value = 0.96319211677255401

$ws.Range("C3").Value = 0.44741899959831699
$ws.Range("D3").Value = 0.55927374949789599
$ws.Range("H3").Value = 0.19773585841432301
$ws.Range("I3").Value = 0.39788977881972598

$ws.Range("C4").Value = 0.0136668244664735
$ws.Range("D4").Value = 0.068334122332367603
$ws.Range("H4").Value = 0.0404900912088582
$ws.Range("I4").Value = 0.20245045604429099

$ws.Range("C5").Value = 0.24968801500461801
$ws.Range("D5").Value = 0.41614669167436402
$ws.Range("H5").Value = 0.96838092421046296
$ws.Range("I5").Value = 0.96838092421046296

$ws.Range("C6").Value = 0.029501308752042699
$ws.Range("D6").Value = 0.073753271880106705
$ws.Range("H6").Value = 0.23873386729183499
$ws.Range("I6").Value = 0.39788977881972598

# D6 now gets called out with a new custom font color (magenta/pink); re-use
# that same highlighted style for D4 too, since its recomputed FDR also now
# stands out.
$ws.Range("D6").Font.Color = 13382655
$ws.Range("D6").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-apply the A-Z sort on HPV.type (data already sorted, but this records
# the sort state the way Excel does after using Sort A to Z on the table).
$rng = $ws.Range("A1:I6")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2"))
$ws.Sort.SetRange($rng)
$ws.Sort.Header = 1
$ws.Sort.Apply()

$ws.Range("G12").Select()
